$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 273
$ws.Range("F4").Value = 9886
$ws.Range("F5").Value = 686
$ws.Range("F6").Value = 178
$ws.Range("F7").Value = 362
$ws.Range("F8").Value = 391
$ws.Range("F9").Value = 440
$ws.Range("F13").Value = 492
$ws.Range("F14").Value = 12501
$ws.Range("F15").Value = 40
$ws.Range("F16").Value = 316
$ws.Range("F19").Value = 255
$ws.Range("F20").Value = 42
$ws.Range("F21").Value = 186
$ws.Range("F22").Value = 131
$ws.Range("F23").Value = 170
$ws.Range("F24").Value = 2747
$ws.Range("F27").Value = 18
$ws.Range("F29").Value = 2162
$ws.Range("F30").Value = 1064
$ws.Range("F31").Value = 4240
$ws.Range("F32").Value = 3756
$ws.Range("F33").Value = 765
$ws.Range("F35").Value = 3074
$ws.Range("F36").Value = 54
$ws.Range("F37").Value = 1351
$ws.Range("F38").Value = 204
$ws.Range("F39").Value = 782
$ws.Range("F40").Value = 38
$ws.Range("F42").Value = 468
$ws.Range("F43").Value = 615
$ws.Range("F45").Value = 151
$ws.Range("F46").Value = 270
$ws.Range("F48").Value = 144
$ws.Range("F49").Value = 159

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 49
$ws.Range("F13").Value = 51

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 60

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 273
$ws.Range("F6").Value = 9886
$ws.Range("F7").Value = 686
$ws.Range("F8").Value = 49
$ws.Range("F9").Value = 178
$ws.Range("F10").Value = 362
$ws.Range("F11").Value = 391
$ws.Range("F12").Value = 440
$ws.Range("F15").Value = 492
$ws.Range("F16").Value = 12501
$ws.Range("F17").Value = 316
$ws.Range("F18").Value = 60
$ws.Range("F19").Value = 255
$ws.Range("F21").Value = 186
$ws.Range("F22").Value = 131
$ws.Range("F23").Value = 170
$ws.Range("F24").Value = 2747
$ws.Range("F28").Value = 2162
$ws.Range("F29").Value = 1064
$ws.Range("F30").Value = 4240
$ws.Range("F31").Value = 3756
$ws.Range("F32").Value = 765
$ws.Range("F34").Value = 3074
$ws.Range("F35").Value = 54
$ws.Range("F36").Value = 1351
$ws.Range("F37").Value = 204
$ws.Range("F38").Value = 782
$ws.Range("F39").Value = 38
$ws.Range("F41").Value = 469
$ws.Range("F43").Value = 615
$ws.Range("F45").Value = 151
$ws.Range("F46").Value = 270
$ws.Range("F48").Value = 144
$ws.Range("F49").Value = 159
